$d = $word.ActiveDocument

$emdash = [char]0x2014
$endash = [char]0x2013

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $t = $r.Text

    if ($t -like "*Maple Labs*Apr 2023*") {
        $r.Text = "Lemon.io " + $emdash + " Full-Stack Developer | Seoul, Korea (Hybrid) | Jan 2024 " + $endash + " Dec 2024"
    }
    elseif ($t -like "*Lead development of onboarding workflows*") {
        $r.Text = "- Shipped onboarding flows, subscription management, and analytics dashboards for the Lemon.io marketplace using Next.js, TypeScript, Tailwind CSS, and GraphQL-backed Node.js services."
    }
    elseif ($t -like "*Partner with product and CX on instrumentation*") {
        $r.Text = "- Partnered with product and CX on instrumentation, ensuring every feature shipped with metrics, feature flags, and rollback plans."
    }
}
